$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended below the existing log (row 86 -> row 87).
# Column A holds a date-like string that must stay plain text (matching
# the existing rows, which are stored as inlineStr/text, not real dates),
# so force a text format before assigning it, then drop the now-unneeded
# number-format override so the cell keeps the sheet's default styling.
$row = 87
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2025/10/10"
$ws.Range("A$row").ClearFormats()

$ws.Range("B$row").Value = "金"
$ws.Range("C$row").Value = 6
$ws.Range("D$row").Value = 201
